$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 / 46: Maker and ThetaToken swap ranking positions ---
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '3.206.68'
$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  -5.05%  '

# --- Price / Volume(1h) refresh for the remaining rows ---
$ws.Range("D2").Value = '69.156.83'
$ws.Range("E2").Value = '  -1.26%  '

$ws.Range("D3").Value = '3.519.06'
$ws.Range("E3").Value = '  -1.71%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.35'
$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.72'
$ws.Range("E6").Value = '  -2.61%  '

$ws.Range("E7").Value = '  -2.48%  '

$ws.Range("D8").Value = '3.512.01'
$ws.Range("E8").Value = '  -1.81%  '

$ws.Range("E10").Value = '  +3.62%  '

$ws.Range("E11").Value = '  -2.90%  '

$ws.Range("E12").Value = '  -3.82%  '

$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.48'
$ws.Range("E14").Value = '  -2.13%  '

$ws.Range("D15").Value = '4.086.20'
$ws.Range("E15").Value = '  -1.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.33'
$ws.Range("E16").Value = '  -3.20%  '

$ws.Range("D17").Value = '3.526.05'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").Value = '69.142.77'
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.52'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("E20").Value = '  -1.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '538.96'
$ws.Range("E21").Value = '  +13.85%  '

$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '20.62'
$ws.Range("E23").Value = '  +7.42%  '

$ws.Range("E24").Value = '  -1.93%  '

$ws.Range("E25").Value = '  +2.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.58'
$ws.Range("E26").Value = '  +6.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("E28").Value = '  -4.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.17'
$ws.Range("E29").Value = '  -2.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.56'
$ws.Range("E30").Value = '  -1.87%  '

$ws.Range("E31").Value = '  -4.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.60'
$ws.Range("E32").Value = '  +4.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.15'
$ws.Range("E33").Value = '  -2.47%  '

$ws.Range("E34").Value = '  -4.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '572.11'
$ws.Range("E35").Value = '  -1.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.19'
$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.08'
$ws.Range("E37").Value = '  +7.59%  '

$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Value = '0.0₃0766'
$ws.Range("E40").Value = '  -4.39%  '

$ws.Range("E41").Value = '  -3.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.10'
$ws.Range("E42").Value = '  -4.48%  '

$ws.Range("E43").Value = '  -4.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.53'
$ws.Range("E44").Value = '  +6.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0440'
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.18'
$ws.Range("E48").Value = '  -3.78%  '

$ws.Range("E49").Value = '  -2.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '136.40'
$ws.Range("E51").Value = '  -0.56%  '
